# Updates cryptos list: refresh Price (D) and Volume(1h) (E) columns for the latest
# snapshot, and apply the USDe/ARBITRUM row swap (rows 45-46) with its own new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 45/46 swapped places: update Coin name + Link first (plain text, no coercion risk) ---
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"

# --- Price (column D): values are stored as plain text (e.g. "1.00", "98.868.21").
#     Toggle to a text format around the write so Excel does not silently coerce
#     number-looking strings ("1.00" -> 1) into numeric cells, then restore the style
#     so no formatting/style change leaks into the saved file. ---
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "98.868.21"
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.341.18"
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "258.23"
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "625.33"
$cell.Style = "Normal"
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.392"
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.881"
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "3.338.73"
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "37.43"
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "98.547.29"
$cell.Style = "Normal"
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "3.959.83"
$cell.Style = "Normal"
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "3.340.24"
$cell.Style = "Normal"
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "3.56"
$cell.Style = "Normal"
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "15.21"
$cell.Style = "Normal"
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "491.03"
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "6.09"
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "9.37"
$cell.Style = "Normal"
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "5.62"
$cell.Style = "Normal"
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "89.34"
$cell.Style = "Normal"
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "11.88"
$cell.Style = "Normal"
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "3.517.26"
$cell.Style = "Normal"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "0.291"
$cell.Style = "Normal"
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.137"
$cell.Style = "Normal"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.995"
$cell.Style = "Normal"
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "28.20"
$cell.Style = "Normal"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "7.27"
$cell.Style = "Normal"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "1.94"
$cell.Style = "Normal"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "499.21"
$cell.Style = "Normal"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "24.88"
$cell.Style = "Normal"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "3.32"
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.780"
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "159.82"
$cell.Style = "Normal"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "45.85"
$cell.Style = "Normal"

# --- Volume(1h) (column E): always padded percent text (e.g. "  +0.93%  "), never
#     numeric-looking, so a direct assignment is safe. ---
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("E3").Value = "  +6.31%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  +6.85%  "
$ws.Range("E6").Value = "  +2.56%  "
$ws.Range("E7").Value = "  +28.04%  "
$ws.Range("E8").Value = "  +2.41%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  +11.80%  "
$ws.Range("E11").Value = "  +6.35%  "
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("E13").Value = "  +10.17%  "
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("E15").Value = "  +3.73%  "
$ws.Range("E16").Value = "  +6.25%  "
$ws.Range("E17").Value = "  +2.01%  "
$ws.Range("E18").Value = "  +6.23%  "
$ws.Range("E19").Value = "  +3.67%  "
$ws.Range("E20").Value = "  +4.67%  "
$ws.Range("E21").Value = "  -5.83%  "
$ws.Range("E22").Value = "  +6.79%  "
$ws.Range("E23").Value = "  +9.68%  "
$ws.Range("E24").Value = "  +6.61%  "
$ws.Range("E25").Value = "  +2.93%  "
$ws.Range("E26").Value = "  +0.95%  "
$ws.Range("E27").Value = "  +2.24%  "
$ws.Range("E28").Value = "  +6.14%  "
$ws.Range("E29").Value = "  +22.78%  "
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("E31").Value = "  +10.72%  "
$ws.Range("E32").Value = "  +11.77%  "
$ws.Range("E33").Value = "  +8.14%  "
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("E35").Value = "  +6.07%  "
$ws.Range("E36").Value = "  -1.41%  "
$ws.Range("E37").Value = "  +0.56%  "
$ws.Range("E38").Value = "  +3.51%  "
$ws.Range("E39").Value = "  +6.67%  "
$ws.Range("E40").Value = "  +5.84%  "
$ws.Range("E41").Value = "  +2.18%  "
$ws.Range("E42").Value = "  +3.47%  "
$ws.Range("E43").Value = "  +4.13%  "
$ws.Range("E44").Value = "  +6.04%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("E46").Value = "  +12.25%  "
$ws.Range("E47").Value = "  -1.43%  "
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("E49").Value = "  +8.80%  "
$ws.Range("E50").Value = "  +2.68%  "
$ws.Range("E51").Value = "  +4.18%  "

Write-Output "cryptos list updated"
